$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on all touched Price/Volume cells so numeric-looking
# strings (e.g. "8.60", "0.0920", "147.10") keep their exact text rendering
# instead of Excel coercing them into numbers and losing formatting.
$cellValues = @{
    "D2" = "62.977.25"
    "E2" = "  +1.52%  "
    "D3" = "2.446.23"
    "E3" = "  +1.07%  "
    "E4" = "  +0.03%  "
    "D5" = "569.84"
    "E5" = "  +1.44%  "
    "D6" = "146.76"
    "E6" = "  +2.28%  "
    "E7" = "  -0.01%  "
    "E8" = "  +0.71%  "
    "D9" = "2.447.31"
    "E9" = "  +1.15%  "
    "E11" = "  +0.36%  "
    "E12" = "  +2.18%  "
    "E13" = "  +0.76%  "
    "E14" = "  +2.62%  "
    "E15" = "  +3.97%  "
    "E16" = "  +1.35%  "
    "D17" = "62.853.59"
    "E17" = "  +1.48%  "
    "D18" = "2.449.63"
    "E18" = "  +1.03%  "
    "D19" = "11.36"
    "E20" = "  +6.17%  "
    "D21" = "324.42"
    "E23" = "  +12.32%  "
    "D24" = "0.999"
    "E24" = "  -0.07%  "
    "D25" = "66.34"
    "E25" = "  -1.49%  "
    "D26" = "620.37"
    "E26" = "  +11.08%  "
    "D27" = "8.60"
    "E27" = "  -1.03%  "
    "D29" = "2.572.49"
    "E29" = "  +1.28%  "
    "D30" = "0.998"
    "E30" = "  -0.23%  "
    "E31" = "  +6.41%  "
    "D32" = "8.27"
    "E32" = "  +0.95%  "
    "D33" = "0.143"
    "E33" = "  -2.54%  "
    "D34" = "1.91"
    "E34" = "  +2.44%  "
    "D35" = "5.09"
    "E35" = "  +7.67%  "
    "E36" = "  -0.57%  "
    "E37" = "  -0.08%  "
    "E38" = "  +0.44%  "
    "E39" = "  -0.99%  "
    "D40" = "18.67"
    "E40" = "  +0.10%  "
    "D41" = "145.87"
    "E41" = "  -4.06%  "
    "D42" = "1.79"
    "E42" = "  -0.62%  "
    "D43" = "2.62"
    "E43" = "  +16.16%  "
    "E44" = "  -0.24%  "
    "D45" = "147.10"
    "E45" = "  -0.12%  "
    "E46" = "  +2.39%  "
    "E47" = "  +2.25%  "
    "D48" = "20.62"
    "E48" = "  +3.64%  "
    "E49" = "  +0.89%  "
    "D50" = "0.0234"
    "E50" = "  +2.36%  "
    "D51" = "0.0920"
    "E51" = "  -0.05%  "
}

foreach ($addr in $cellValues.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $cellValues[$addr]
}
